$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix product-name typo: "Squalene" -> "Squalane" (rows 14-16, Sr. No. 13-15)
$ws.Range("B14:B16").Value = "Squalane Glow Moisturizer"

# Clear the (redundant/no-op) explicit style on the "BUY2649" offer cells in
# column C so they fall back to the sheet's default style, matching the
# re-saved workbook's pruned style table.
$buyRows = @(2,5,8,11,14,17,20,23,26,29,32,35,38,41)
foreach ($r in $buyRows) {
    $ws.Cells.Item($r, 3).Style = "Normal"
}

# Restore the cursor/selection state left by the editing session.
$ws.Range("B17").Select() | Out-Null
